$d = $word.ActiveDocument

$pairs = @(
    @("65×63=4095", "87×94=8178"),
    @("65×99=6435", "24×79=1896"),
    @("27×92=2484", "49×89=4361"),
    @("64×72=4608", "11×35=385"),
    @("16×97=1552", "46×36=1656"),
    @("60×29=1740", "52×15=780"),
    @("96×73=7008", "83×39=3237"),
    @("13×21=273", "18×11=198"),
    @("43×52=2236", "71×58=4118"),
    @("42×12=504", "57×39=2223"),
    @("75×41=3075", "16×62=992"),
    @("90×95=8550", "13×54=702"),
    @("75×72=5400", "47×33=1551"),
    @("71×94=6674", "60×57=3420"),
    @("36×88=3168", "93×53=4929"),
    @("16×16=256", "23×51=1173"),
    @("40×46=1840", "30×22=660"),
    @("92×79=7268", "59×72=4248"),
    @("26×61=1586", "60×74=4440"),
    @("86×77=6622", "46×12=552"),
    @("83×12=996", "83×13=1079"),
    @("22×49=1078", "52×93=4836"),
    @("66×15=990", "75×88=6600"),
    @("23×50=1150", "64×85=5440"),
    @("88×71=6248", "54×20=1080")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
